# "feat: new bot created"
# The workbook (an "Open OMs" export, one row per open order) gains a new
# "Status" style column and an "Unnamed: 0" index column, matching the
# output of a bot that re-exported the sheet from pandas and then
# annotated each order with its current status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) New column A header ("Unnamed: 0" - the pandas index column),
#    cloning B1's header style (bold, bordered, centered) so no new
#    style slot is introduced.
# ---------------------------------------------------------------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "OM"
$ws.Range("C1").Value = "Status"

# ---------------------------------------------------------------
# 2) The old index column (A2:A59) drops its special formatting -
#    it is now a plain pandas index column with no explicit style.
# ---------------------------------------------------------------
$ws.Range("A2:A59").ClearFormats() | Out-Null

# ---------------------------------------------------------------
# 3) The OM numbers in column B are now centered instead of left
#    aligned.
# ---------------------------------------------------------------
$ws.Range("B2:B59").HorizontalAlignment = -4108

# ---------------------------------------------------------------
# 4) New "Status" column C: each open order that has already been
#    processed is flagged "Encerrado!" (closed), the one still
#    awaiting action is flagged "Ordem pendente!" (pending order),
#    and orders not yet checked are left blank (but still styled /
#    centered, ready to be filled in later).
# ---------------------------------------------------------------
$statusMap = @{}
$statusMap[2]  = "Encerrado!"
$statusMap[3]  = "Encerrado!"
$statusMap[4]  = "Encerrado!"
$statusMap[5]  = "Encerrado!"
$statusMap[6]  = "Encerrado!"
$statusMap[7]  = "Encerrado!"
$statusMap[8]  = "Encerrado!"
$statusMap[9]  = "Encerrado!"
$statusMap[10] = "Encerrado!"
$statusMap[11] = "Encerrado!"
$statusMap[12] = "Encerrado!"
$statusMap[13] = "Encerrado!"
$statusMap[14] = "Encerrado!"
$statusMap[15] = "Encerrado!"
$statusMap[16] = "Encerrado!"
$statusMap[17] = "Encerrado!"
$statusMap[18] = "Encerrado!"
$statusMap[19] = "Encerrado!"
$statusMap[20] = "Encerrado!"
$statusMap[21] = "Ordem pendente!"
$statusMap[22] = "Encerrado!"
$statusMap[23] = "Encerrado!"
$statusMap[24] = "Encerrado!"
$statusMap[25] = "Encerrado!"

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = $statusMap[$r]
}

$ws.Range("C2:C59").HorizontalAlignment = -4108

# ---------------------------------------------------------------
# 5) Column C sizes itself to fit the new status text.
# ---------------------------------------------------------------
$ws.Columns.Item(3).AutoFit() | Out-Null

# ---------------------------------------------------------------
# 6) Misc bookkeeping matching a normal "re-saved after review" pass:
#    the selection ends up on the next free row, and the page
#    margins are reset to Excel's defaults.
# ---------------------------------------------------------------
$ws.Range("E76").Select() | Out-Null

$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
